$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1

$ws.Range("A3").Value = "10uH inductor"
$ws.Range("B3").Value = "587-2886-1-ND"
$ws.Range("C3").Value = 0.29
$ws.Range("D3").Value = 1

$ws.Range("D1").Value = "Needed"

$ws.Range("A4").Value = "47uF ceramic"
$ws.Range("D4").Value = 1

$ws.Range("A5").Value = "4.7uF ceramic"
$ws.Range("D5").Value = 1

$ws.Range("E12").Select()
